$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C6").ClearFormats()

$ws.Range("C2").Value = 400
$ws.Range("C3").Value = 800
$ws.Range("C4").Value = 8900
$ws.Range("C5").Value = 56
$ws.Range("C6").Value = 78
